$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Extend the assumptions table from year 2037 (row 18) through year 2060 (row 41),
# repeating the same pattern used for the existing rows: binomia=1, alternativa=1,
# p_transicao=FALSE, demanda_g=TRUE.
$startRow = 19
$startYear = 2038
$endYear = 2060

for ($i = 0; $i -le ($endYear - $startYear); $i++) {
    $row = $startRow + $i
    $year = $startYear + $i

    $ws.Cells.Item($row, 1).Value = $year
    $ws.Cells.Item($row, 2).Value = 1
    $ws.Cells.Item($row, 3).Value = 1
    $ws.Cells.Item($row, 4).Formula = "=FALSE"
    $ws.Cells.Item($row, 5).Formula = "=TRUE"
}
